$p = $ppt.ActivePresentation

# Remove the last slide (sldId 262 / slide6.xml) from the deck.
$p.Slides.Item(6).Delete()

# The "datetimeFigureOut" date placeholder text on the slide master and every
# slide layout was re-cached from 2024/3/3 to 2024/2/13.
$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = "2024/2/13"
    }
}

$cls = $m.CustomLayouts
for ($j = 1; $j -le $cls.Count; $j++) {
    $cl = $cls.Item($j)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = "2024/2/13"
        }
    }
}
